# Generate Report for Handoff
# Updates the row for "b1ca587b-6455-4767-838f-95f283b008f0" across the
# Overview, zh-cn and de-de sheets to reflect that the file is now
# "Ready for handoff" with refreshed handoff timestamps.

$wb = $excel.ActiveWorkbook

# --- Overview sheet (row 3 = b1ca587b-6455-4767-838f-95f283b008f0.md) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-14-18 00:14:44"

# --- zh-cn sheet (row 3 = b1ca587b-6455-4767-838f-95f283b008f0.md) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-18 00:14:41"

# --- de-de sheet (row 3 = b1ca587b-6455-4767-838f-95f283b008f0.md) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-18 00:14:44"
